$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (a new "2007 -> 2008" forecast
# record), shifting the existing year rows (previously rows 2-18) down to
# rows 3-19.
$ws.Rows.Item(2).Insert()

# The native row-insert copies row 1's (header) formatting across the whole
# new row 2. Re-apply the correct data-row formatting (date style on column
# A, default style on B:E) by copying the format from row 3, which already
# carries the right per-column styles.
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Full corrected data set: row, date-serial (A), y_0 (B), y_0_forecast (C),
# y_1 (D), y_1_forecast (E).
$data = @(
    @(2, 39400, 2007, 7.226520411029047, 2008, 8.987952903094421),
    @(3, 39765, 2008, 4.268860212333636, 2009, 6.77211531652997),
    @(4, 40130, 2009, -7.266312015249776, 2010, -1.420242831007679),
    @(5, 40494, 2010, 6.958243460951929, 2011, 8.31992099499319),
    @(6, 40862, 2011, 9.469137444079934, 2012, 8.571528775834981),
    @(7, 41228, 2012, 3.358206407534947, 2013, 5.745831525574441),
    @(8, 41592, 2013, 0.3081076735359067, 2014, 2.847379875994704),
    @(9, 41957, 2014, 3.901355411819707, 2015, 3.690459963535009),
    @(10, 42321, 2015, 5.331683351557981, 2016, 4.636575318346536),
    @(11, 42689, 2016, 3.254758369308375, 2017, 2.76788332063731),
    @(12, 43053, 2017, 5.246209615995667, 2018, 4.659862065074982),
    @(13, 43418, 2018, 4.86255966374296, 2019, 4.636196713604357),
    @(14, 43783, 2019, 2.764740011159428, 2020, 2.471557257221946),
    @(15, 44159, 2020, -7.260793671746435, 2021, -1.387795042833839),
    @(16, 44525, 2021, 4.097586525396268, 2022, 2.450242954096926),
    @(17, 44890, 2022, 7.824284864703746, 2023, 4.834990656989402),
    @(18, 45254, 2023, -1.24502235313334, 2024, -2.488220481262082),
    @(19, 45618, 2024, -1.735114423676209, 2025, 1.056286187957367)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
